$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4, col A was recorded as text; correct it to a numeric value.
$ws.Range("A4").Value = 79174445

# Append new redemption row (row 5): phone 79174445 redeems 20 points.
$ws.Range("A5").Value = "'79174445"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = "2025-08-18T08:51:56"
